$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B31: convert from text "4" to a real number 4
$ws.Range("B31").Value = 4

# Add new row 32 with the new annotation data
$ws.Range("A32").Value = "Ying Tang"

# B32 must stay text "5" (not a number) - use a leading apostrophe to force
# text entry, then reset the style so no extra formatting/quote-prefix is
# left behind on the cell.
$ws.Range("B32").Value = "'5"
$ws.Range("B32").Style = "Normal"

$ws.Range("C32").Value = "thank ,detailed and insightful feedback"
$ws.Range("D32").Value = "APC"
$ws.Range("E32").Value = "OTH"
$ws.Range("F32").Value = "7f314748-ac5a-4a11-8786-6125314f9d6d"
$ws.Range("G32").Value = "Sy2ogebAW_annotated.xlsx"
$ws.Range("H32").Value = "We would like to thank all reviewers for their detailed and insightful feedback."
